# dwr-offices.xlsx: correct the Division 1 Office URL and drop the now
# broken hyperlink on the Data sheet (the old water.state.co.us page has
# moved to dwr.colorado.gov).

$wb = $excel.ActiveWorkbook

$data = $wb.Worksheets.Item("Data")

# Cell B2 shows the URL as its own display text; update it to the new
# DWR site and remove the now-stale hyperlink attached to the cell while
# keeping the cell's existing (Hyperlink) style untouched.
$data.Range("B2").Value = "https://dwr.colorado.gov/division-offices/division-1-office"
$data.Range("B2").Hyperlinks.Delete()
